$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.238.91"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.353.12"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'189.96"
$ws.Range("E5").Value = "  +5.60%  "
$ws.Range("D6").Value = "'557.77"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.345.71"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'46.52"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "3.890.47"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "'8.56"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "'594.29"
$ws.Range("E16").Value = "  -6.15%  "
$ws.Range("D17").Value = "66.312.47"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "3.344.66"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("D19").Value = "'17.95"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "'11.04"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").Value = "'0.901"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'18.37"
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("D24").Value = "'5.02"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "'99.15"
$ws.Range("E25").Value = "  -6.50%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").Value = "'6.06"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("D29").Value = "'9.49"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'8.53"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").Value = "'30.83"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").Value = "'6.71"
$ws.Range("E32").Value = "  +6.24%  "
$ws.Range("D33").Value = "'3.83"
$ws.Range("E33").Value = "  -5.97%  "
$ws.Range("D34").Value = "'583.17"
$ws.Range("E34").Value = "  +6.53%  "
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "3.772.17"
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "'55.98"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").Value = "'34.40"
$ws.Range("E40").Value = "  +7.82%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.127"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0699"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").Value = "'3.17"
$ws.Range("E44").Value = "  -6.90%  "
$ws.Range("D45").Value = "'3.39"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0416"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("B48").Value = "CoreDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D48").Value = "'3.12"
$ws.Range("E48").Value = "  -15.88%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("E51").Value = "  +0.32%  "
